$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 149.6
$ws.Range("I11").Value = 149.6
$ws.Range("K11").Value = 149.6
$ws.Range("M11").Value = -9.599999999999994
$ws.Range("H17").Value = 48896.523
$ws.Range("J17").Value = 48896.523
$ws.Range("L17").Value = 146689.569
$ws.Range("N17").Value = -147025.569
$ws.Range("H33").Value = 294.7143
$ws.Range("I33").Value = 294.7143
$ws.Range("K33").Value = 294.7143
$ws.Range("M33").Value = -65.71429999999998
$ws.Range("H40").Value = 62502460
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 71431096
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 71431096
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -71431446
$ws.Range("H98").Value = 1721.4359
$ws.Range("I98").Value = 1392.3611
$ws.Range("K98").Value = 1392.3611
$ws.Range("M98").Value = 105.6388999999999
$ws.Range("H122").Value = 1721.4359
$ws.Range("I122").Value = 1392.3611
$ws.Range("K122").Value = 4177.0833
$ws.Range("M122").Value = -1727.0833
$ws.Range("H132").Value = 1275.1086
$ws.Range("I132").Value = 1287.9111
$ws.Range("K132").Value = 3863.7333
$ws.Range("M132").Value = -1333.7333
$ws.Range("H138").Value = 3325.8357
$ws.Range("I138").Value = 2245.8965
$ws.Range("J138").Value = 4037.6135
$ws.Range("K138").Value = 6737.689499999999
$ws.Range("L138").Value = 12112.8405
$ws.Range("M138").Value = -1597.689499999999
$ws.Range("N138").Value = -22392.8405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1973.6316
$ws.Range("I2").Value = 1933.8667
$ws.Range("K2").Value = 1933.8667
$ws.Range("M2").Value = -1820.8667
$ws.Range("H61").Value = 3446488
$ws.Range("I61").Value = 3574983.5
$ws.Range("K61").Value = 3574983.5
$ws.Range("M61").Value = -3574771.5
$ws.Range("H74").Value = 2860
$ws.Range("I74").Value = 2825.2942
$ws.Range("K74").Value = 2825.2942
$ws.Range("M74").Value = -1951.2942
$ws.Range("H77").Value = 2860
$ws.Range("I77").Value = 2825.2942
$ws.Range("K77").Value = 14126.471
$ws.Range("M77").Value = -9758.471
$ws.Range("H116").Value = 1973.6316
$ws.Range("I116").Value = 1933.8667
$ws.Range("K116").Value = 1933.8667
$ws.Range("M116").Value = 360.1333
$ws.Range("H132").Value = 1286647
$ws.Range("I132").Value = 4473.9395
$ws.Range("K132").Value = 13421.8185
$ws.Range("M132").Value = -10891.8185
$ws.Range("H136").Value = 3446488
$ws.Range("I136").Value = 3574983.5
$ws.Range("K136").Value = 10724950.5
$ws.Range("M136").Value = -10722400.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1973.6316
$ws.Range("I3").Value = 1933.8667
$ws.Range("K3").Value = 1933.8667
$ws.Range("M3").Value = -1819.8667
$ws.Range("H80").Value = 1070.4231
$ws.Range("I80").Value = 1304.2
$ws.Range("J80").Value = 924.3125
$ws.Range("K80").Value = 1304.2
$ws.Range("L80").Value = 924.3125
$ws.Range("M80").Value = -306.2
$ws.Range("N80").Value = -2920.3125
$ws.Range("H83").Value = 1070.4231
$ws.Range("I83").Value = 1304.2
$ws.Range("J83").Value = 924.3125
$ws.Range("K83").Value = 6521
$ws.Range("L83").Value = 4621.5625
$ws.Range("M83").Value = -1529
$ws.Range("N83").Value = -14605.5625
$ws.Range("H86").Value = 4021.5
$ws.Range("I86").Value = 3416
$ws.Range("K86").Value = 3416
$ws.Range("M86").Value = -2293
$ws.Range("H89").Value = 4021.5
$ws.Range("I89").Value = 3416
$ws.Range("K89").Value = 17080
$ws.Range("M89").Value = -11464

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 19173.4
$ws.Range("I51").Value = 19173.4
$ws.Range("K51").Value = 19173.4
$ws.Range("M51").Value = -18437.4
$ws.Range("H61").Value = 19173.4
$ws.Range("I61").Value = 19173.4
$ws.Range("K61").Value = 19173.4
$ws.Range("M61").Value = -18825.4
$ws.Range("H105").Value = 1312.8889
$ws.Range("I105").Value = 961.2
$ws.Range("J105").Value = 1752.5
$ws.Range("K105").Value = 961.2
$ws.Range("L105").Value = 1752.5
$ws.Range("M105").Value = 785.8
$ws.Range("N105").Value = -5246.5
$ws.Range("H129").Value = 99999
$ws.Range("J129").Value = 99999
$ws.Range("L129").Value = 99999
$ws.Range("N129").Value = -109999
$ws.Range("H141").Value = 235297
$ws.Range("J141").Value = 235297
$ws.Range("L141").Value = 235297
$ws.Range("N141").Value = -245657

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 14964.6
$ws.Range("I123").Value = 11495
$ws.Range("K123").Value = 34485
$ws.Range("M123").Value = -32035

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5734.9355
$ws.Range("I122").Value = 5161.4287
$ws.Range("J122").Value = 6939.3
$ws.Range("K122").Value = 15484.2861
$ws.Range("L122").Value = 20817.9
$ws.Range("M122").Value = -13034.2861
$ws.Range("N122").Value = -25717.9
$ws.Range("H132").Value = 2704244.8
$ws.Range("I132").Value = 1404.2727
$ws.Range("K132").Value = 4212.8181
$ws.Range("M132").Value = -1682.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8191.8
$ws.Range("J22").Value = 657
$ws.Range("L22").Value = 657
$ws.Range("N22").Value = -1247
$ws.Range("H27").Value = 8191.8
$ws.Range("J27").Value = 657
$ws.Range("L27").Value = 657
$ws.Range("N27").Value = -871
$ws.Range("H43").Value = 19142.715
$ws.Range("J43").Value = 89999
$ws.Range("L43").Value = 89999
$ws.Range("N43").Value = -90385
$ws.Range("H46").Value = 1212.5
$ws.Range("I46").Value = 1141.6666
$ws.Range("K46").Value = 1141.6666
$ws.Range("M46").Value = -953.6666
$ws.Range("H61").Value = 4606.316
$ws.Range("I61").Value = 3911.077
$ws.Range("K61").Value = 3911.077
$ws.Range("M61").Value = -3709.077
$ws.Range("H68").Value = 4169646.2
$ws.Range("I68").Value = 6946742
$ws.Range("K68").Value = 6946742
$ws.Range("M68").Value = -6945993
$ws.Range("H71").Value = 4169646.2
$ws.Range("I71").Value = 6946742
$ws.Range("K71").Value = 34733710
$ws.Range("M71").Value = -34729966
$ws.Range("H113").Value = 4606.316
$ws.Range("I113").Value = 3911.077
$ws.Range("K113").Value = 3911.077
$ws.Range("M113").Value = -1741.077
$ws.Range("H132").Value = 4314.3125
$ws.Range("I132").Value = 2485.5557
$ws.Range("J132").Value = 6665.5713
$ws.Range("K132").Value = 7456.6671
$ws.Range("L132").Value = 19996.7139
$ws.Range("M132").Value = -4926.6671
$ws.Range("N132").Value = -25056.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2802.0715
$ws.Range("I126").Value = 3135.3684
$ws.Range("K126").Value = 9406.1052
$ws.Range("M126").Value = -6936.1052
